# INDIANA_2016.xlsx cleanup:
#  1. Rename header row (A1:D1) to short machine-friendly column names.
#  2. Title-case the Spanish joining words ("de", "del", "la", "las", "los",
#     "el", "y") inside the state/municipality name columns (they were
#     previously lower-cased except when first word of the cell).
#  3. Fix two floating point rounding artifacts in column D.
#  4. Drop the trailing footnote/metadata rows (1323-1327) and shrink the
#     used range back down to A1:D1321.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header renames -------------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Title-case the Spanish connector words ------------------------------
# Space-padded whole-word substitution so the leading word of a cell
# (e.g. "El Llano", "La Paz") is never touched, and mid-word substrings
# (e.g. "del" inside another token) can't match either.
$joiners = @("de", "del", "la", "las", "los", "el", "y")
foreach ($w in $joiners) {
    $find = " " + $w + " "
    $firstLetter = $w.Substring(0,1).ToUpper()
    $rest = $w.Substring(1)
    $replace = " " + $firstLetter + $rest + " "
    $ws.Cells.Replace($find, $replace, [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart, [Microsoft.Office.Interop.Excel.XlSearchOrder]::xlByRows, $true, $false, $true)
}

# 3. Floating point precision fixes ---------------------------------------
$ws.Range("D149").Value2 = 0.009370494954348873
$ws.Range("D1083").Value2 = 0.009290405253884353

# 4. Drop the trailing metadata rows --------------------------------------
$ws.Range("A1323:A1327").EntireRow.Delete()
